# Mvmnt File degree/radian option added. Pose extraction for mult markers added
#
# This adds a brand-new worksheet "many_markers" right after "Main" (i.e. as
# the 2nd tab), shifting every other sheet one position to the right. The new
# sheet mirrors the same "Parameter / Info / Additional_Info" layout used by
# the other simulation-config sheets, but documents a marker set with many
# more entries (multiple dictionaries / multiple marker poses).

$wb = $excel.ActiveWorkbook

# Use "mult_dict" as a structural template (same column widths/styles/layout)
# and copy it right after "Main" - this becomes our new "many_markers" tab.
$template = $wb.Worksheets.Item("mult_dict")
$template.Copy($null, $wb.Worksheets.Item("Main"))

$ws = $wb.Worksheets.Item(2)
$ws.Name = "many_markers"

# Wipe everything below the shared header row so stale template values can't
# leak through, then (re)populate the full A1:C16 block for this sheet.
$ws.Range("A2:C16").ClearContents()

$ws.Range("A1").Value = "Parameter"
$ws.Range("B1").Value = "Info"
$ws.Range("C1").Value = "Additional_Info"

$ws.Range("A2").Value = "movement_file"
$ws.Range("B2").Value = "many_markers.txt"

$ws.Range("A3").Value = "video_file"
$ws.Range("B3").Value = 1

$ws.Range("A4").Value = "gz_pose_file"
$ws.Range("B4").Value = 1

$ws.Range("A5").Value = "vid_pose_file"
$ws.Range("B5").Value = 1

$ws.Range("A6").Value = "cameras"
$ws.Range("B6").Value = "KahnPhone_new.sdf"
$ws.Range("C6").Value = "1,1,1,0,0,0"

$ws.Range("A7").Value = "markers"
$ws.Range("B7").Value = "DICT_4X4_50_s100_id0.sdf"
$ws.Range("C7").Value = "0,0,0,0,0,0"

$ws.Range("B8").Value = "DICT_4X4_50_s500_id1.sdf"
$ws.Range("C8").Value = "0.7,0.2,0,0,0,0.7"

$ws.Range("B9").Value = "DICT_4X4_50_s1000_id3.sdf"
$ws.Range("C9").Value = "0,-2.3,0.8,-0.5,0.3,0.8"

$ws.Range("B10").Value = "DICT_4X4_50_s100_id9.sdf"
$ws.Range("C10").Value = "-0.1,-0.2,0.1,-0.2,0.2,0"

$ws.Range("B11").Value = "DICT_4X4_50_s100_id5.sdf"
$ws.Range("C11").Value = "1,1,0,0,0,0"

$ws.Range("B12").Value = "DICT_5X5_50_s100_id0.sdf"
$ws.Range("C12").Value = "0,0.5,1,0,0,0"

$ws.Range("B13").Value = "DICT_5X5_50_s200_id8.sdf"
$ws.Range("C13").Value = "-1.2,0.5,0,0,-0.53,2.5"

$ws.Range("B14").Value = "DICT_5X5_50_s100_id4.sdf"
$ws.Range("C14").Value = "-0.25,0,0.1,0,0.52,0"

$ws.Range("A15").Value = "lights"
$ws.Range("A16").Value = "models"

# Make the freshly-added sheet the active tab, matching the authored file's
# selection (cursor sitting on B2, the movement-file cell).
$ws.Activate() | Out-Null
$ws.Range("B2").Select() | Out-Null
